$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.294.96"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").Value = "1.864.97"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.84"
$ws.Range("E5").Value = "  -0.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4980"
$ws.Range("E7").Value = "  -3.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3981"
$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09972"
$ws.Range("E9").Value = "  +27.51%  "

$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.475"
$ws.Range("E12").Value = "  +2.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.85"
$ws.Range("E13").Value = "  +2.02%  "

$ws.Range("D14").Value = "1.860.98"
$ws.Range("E14").Value = "  +3.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.002"
$ws.Range("E15").Value = "  -0.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.403"
$ws.Range("E16").Value = "  +1.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001143"
$ws.Range("E17").Value = "  +5.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.52"
$ws.Range("E18").Value = "  +1.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06648"
$ws.Range("E19").Value = "  +1.16%  "

$ws.Range("E21").Value = "  +0.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.074"

$ws.Range("D23").Value = "28.367.43"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.251"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.19"
$ws.Range("E26").Value = "  +3.36%  "

$ws.Range("D27").Value = "2.075.64"
$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.481"
$ws.Range("E28").Value = "  +2.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.44"
$ws.Range("E29").Value = "  -2.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.74"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("E31").Value = "  -3.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").Value = "  -1.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.628"
$ws.Range("E33").Value = "  +0.95%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.599"
$ws.Range("E34").Value = "  -1.70%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06807"
$ws.Range("E35").Value = "  -5.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.185"
$ws.Range("E36").Value = "  +0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02381"
$ws.Range("E37").Value = "  +0.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2164"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.020"
$ws.Range("E39").Value = "  -0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.47"

$ws.Range("E41").Value = "  +1.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.179"
$ws.Range("E42").Value = "  +1.96%  "

$ws.Range("E43").Value = "  -0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.43"
$ws.Range("E44").Value = "  +1.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5993"
$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.280"
$ws.Range("E46").Value = "  -2.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.671"
$ws.Range("E47").Value = "  -1.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.67"
$ws.Range("E48").Value = "  -0.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.981"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.188"
$ws.Range("E50").Value = "  -2.30%  "

$ws.Range("E51").Value = "  +4.09%  "
